$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '93.638.17'
$ws.Range('E2').Value = '  +3.83%  '
$ws.Range('D3').Value = '3.134.05'
$ws.Range('E3').Value = '  -0.38%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = '243.54'
$ws.Range('E5').Value = '  +2.84%  '
$ws.Range('D6').Value = '618.93'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').Value = '0.411'
$ws.Range('E8').Value = '  +11.55%  '
$ws.Range('D9').Value = '0.998'
$ws.Range('E9').Value = '  -0.24%  '
$ws.Range('D10').Value = '3.129.52'
$ws.Range('E10').Value = '  -0.48%  '
$ws.Range('D11').Value = '0.744'
$ws.Range('E11').Value = '  +1.47%  '
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').Value = '''0.0000257'
$ws.Range('E13').Value = '  +4.81%  '
$ws.Range('E14').Value = '  -0.42%  '
$ws.Range('D15').Value = '93.100.68'
$ws.Range('E15').Value = '  +2.97%  '
$ws.Range('D16').Value = '''5.50'
$ws.Range('E16').Value = '  -0.21%  '
$ws.Range('D17').Value = '3.717.49'
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('D18').Value = '3.120.74'
$ws.Range('E18').Value = '  -1.30%  '
$ws.Range('D19').Value = '3.78'
$ws.Range('E19').Value = '  +2.91%  '
$ws.Range('D20').Value = '14.91'
$ws.Range('E20').Value = '  -1.20%  '
$ws.Range('D21').Value = '''0.0000210'
$ws.Range('E21').Value = '  +3.73%  '
$ws.Range('D22').Value = '5.85'
$ws.Range('E22').Value = '  +0.63%  '
$ws.Range('E23').Value = '  +5.04%  '
$ws.Range('D24').Value = '452.39'
$ws.Range('E24').Value = '  +3.20%  '
$ws.Range('E25').Value = '  +2.85%  '
$ws.Range('D26').Value = '87.95'
$ws.Range('E26').Value = '  -0.70%  '
$ws.Range('D27').Value = '11.94'
$ws.Range('E27').Value = '  +1.13%  '
$ws.Range('D28').Value = '3.294.73'
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('E30').Value = '  +8.44%  '
$ws.Range('D31').Value = '''0.170'
$ws.Range('E31').Value = '  +0.30%  '
$ws.Range('E32').Value = '  -0.28%  '
$ws.Range('D33').Value = '''9.30'
$ws.Range('E33').Value = '  +0.16%  '
$ws.Range('D34').Value = '0.996'
$ws.Range('E34').Value = '  -0.52%  '
$ws.Range('E35').Value = '  +6.18%  '
$ws.Range('D36').Value = '0.163'
$ws.Range('E36').Value = '  -3.04%  '
$ws.Range('D37').Value = '26.42'
$ws.Range('E37').Value = '  +1.58%  '
$ws.Range('E38').Value = '  -0.76%  '
$ws.Range('D39').Value = '3.95'
$ws.Range('E39').Value = '  +4.78%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').Value = '1.32'
$ws.Range('E40').Value = '  -1.68%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').Value = '484.83'
$ws.Range('E41').Value = '  -3.44%  '
$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D42').Value = '0.438'
$ws.Range('E42').Value = '  -1.73%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').Value = '3.47'
$ws.Range('E43').Value = '  +1.49%  '
$ws.Range('E44').Value = '  +4.60%  '
$ws.Range('D46').Value = '162.53'
$ws.Range('E46').Value = '  +2.61%  '
$ws.Range('D47').Value = '1.96'
$ws.Range('E47').Value = '  +2.64%  '
$ws.Range('E48').Value = '  -1.19%  '
$ws.Range('E49').Value = '  +3.16%  '
$ws.Range('D50').Value = '0.0337'
$ws.Range('E50').Value = '  +4.90%  '
$ws.Range('D51').Value = '''4.50'
$ws.Range('E51').Value = '  +2.18%  '
